# Add a new "Keras and TensorFlow" bullet right after the
# "Linear and logistic regression with a single neuron" item
# (the last paragraph of the document), matching the same list
# level/formatting, and reproducing the Word spell-check markup
# (<w:proofErr>) around the unrecognised word "Keras".

$d = $word.ActiveDocument

# Collapse a range to the end of the document content and insert the
# new paragraph there as raw OOXML so we can control the exact run
# layout (including the proofErr spell-check markers) while reusing
# the same paragraph/run formatting (Arial font, en-US language,
# ListParagraph style, list level 0 of numId 4).
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr>' +
    '<w:pStyle w:val="ListParagraph"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr>' +
  '</w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' +
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t>Keras</w:t>' +
  '</w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' +
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t xml:space="preserve"> and TensorFlow</w:t>' +
  '</w:r>' +
  '</w:p>'

$null = $insertionPoint.InsertXML($newParagraphXml)
